# edit.ps1
# Applies the TR.docx content/style changes described by the target diff:
#  1. Renames three bookmarks from their generated hash-like names to slugs.
#  2. Reformats the three "organization:" quote paragraphs from BlockText
#     (soft-wrapped, curly-quoted) paragraphs into SourceCode paragraphs made
#     of VerbatimChar runs separated by explicit text-wrapping line breaks,
#     with straight quotes and indentation preserved from the source.
#  3. Replaces "Example Org" with "Project" throughout (also collapses the
#     now-contiguous runs in those paragraphs into single runs, matching
#     Word's normal Find&Replace behavior).

$d = $word.ActiveDocument

function Rename-Bookmark {
    param($doc, $oldName, $newName)
    $bm = $doc.Bookmarks.Item($oldName)
    $rng = $bm.Range
    $bm.Delete()
    $doc.Bookmarks.Add($newName, $rng) | Out-Null
}

function Set-VerbatimParagraph {
    param($doc, $paraIndex, $lines)
    $p = $doc.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = ""
    $p.Style = "SourceCode"

    $nl = [char]11
    $full = $lines -join $nl

    $rng2 = $doc.Paragraphs.Item($paraIndex).Range
    $rng2.MoveEnd(1, -1) | Out-Null
    $rng2.InsertAfter($full)

    $pos = $rng2.Start
    foreach ($line in $lines) {
        $len = $line.Length
        $lr = $doc.Range($pos, $pos + $len)
        $lr.Style = "VerbatimChar"
        $pos = $pos + $len + 1
    }
}

$lines5 = @(
  '"The organization:',
  '     a.    Provides effective notice to the public and to individuals regarding:',
  '(i) its activities that',
  '           impact privacy, including its collection, use, sharing, safeguarding,',
  'maintenance, and disposal',
  '           of personally identifiable information (PII); (ii) authority for collecting',
  'PII; (iii) the choices, if',
  '           any, individuals may have regarding how the organization uses PII and',
  'the consequences of',
  '           exercising or not exercising those choices; and (iv) the ability to',
  'access and have PII amended',
  '           or corrected if necessary;',
  '     b.    Describes: (i) the PII the organization collects and the purpose(s)',
  'for which it collects that',
  '           information; (ii) how the organization uses PII internally; (iii) whether',
  'the organization shares',
  '           PII with external entities, the categories of those entities, and the',
  'purposes for such sharing;',
  '           (iv) whether individuals have the ability to consent to specific uses',
  'or sharing of PII and how',
  '           to exercise any such consent; (v) how individuals may obtain access',
  'to PII; and (vi) how the',
  '           PII will be protected; and',
  '     c.    Revises its public notices to reflect changes in practice or policy',
  'that affect PII or changes in',
  '           its activities that impact privacy, before or as soon as practicable',
  'after the change."'
)

$lines9 = @(
  '"The organization:',
  '     a.    Publishes System of Records Notices (SORNs) in the Federal Register,',
  'subject to required',
  '           oversight processes, for systems containing personally identifiable',
  'information (PII);',
  '     b.    Keeps SORNs current; and',
  '     c.    Includes Privacy Act Statements on its forms that collect PII, or on',
  'separate forms that can be',
  '           retained by individuals, to provide additional formal notice to individuals',
  'from whom the',
  '           information is being collected."'
)

$lines13 = @(
  '"The organization:',
  '     a.    Ensures that the public has access to information about its privacy activities and is able to',
  '           communicate with its Senior Agency Official for Privacy (SAOP)/Chief Privacy Officer',
  '           (CPO); and',
  '     b.    Ensures that its privacy practices are publicly available through organizational websites or',
  '           otherwise."'
)


# --- 1. Bookmark renames ---
Rename-Bookmark $d "X4569eb5fcd13585a93b272e0a5ded52ee0dd8ee" "reusable-component-library-system-security-plan"
Rename-Bookmark $d "X00439c23c0d5b170b83f352bd0c5977894fb2f5" "tr-2-system-of-records-notices-and-privacy-act-statements"
Rename-Bookmark $d "X0b7aa0d460fb030f2067ca8b3740803deb6d8a7" "tr-3-dissemination-of-privacy-program-information"

# --- 2. Reformat the three quoted-requirement paragraphs ---
Set-VerbatimParagraph $d 5 $lines5
Set-VerbatimParagraph $d 9 $lines9
Set-VerbatimParagraph $d 13 $lines13

# --- 3. "Example Org" -> "Project" everywhere ---
$d.Content.Find.Execute("Example Org", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Project", 2) | Out-Null

Write-Output "edit complete"
